# Apply "feat: add 2022-Q3 data" to the workbook.
#
# Net effect (by sheet NAME, not by package file position):
#   - "总计" (summary) sheet gets a new top data row for "2022-Q3" inserted
#     above the existing quarters, which all shift down one row but keep
#     their own values.
#   - A brand-new sheet "2022-Q3" is inserted right after "总计", holding
#     the fund holdings table for the new quarter.
#   - The existing "2022-Q2", "2022-Q1" and "2021-Q4" sheets are left
#     completely untouched (same name, same data, just pushed one tab to
#     the right to make room for the new quarter).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: make room for 2022-Q3 at the top of the table
#    and re-write the (small, fully known) A2:D5 block with final values.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Give row 5 the same formatting as row 4 (thin border / bold / centered
# index cell in column A) before filling it with data, so the newly
# "grown" table row looks like the others instead of picking up default
# formatting.
$zj.Range("A4:D4").Copy()
$zj.Range("A5:D5").PasteSpecial(-4122)   # xlPasteFormats

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 0.2

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q2"
$zj.Range("C3").Value = 3
$zj.Range("D3").Value = 0.28

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q1"
$zj.Range("C4").Value = 3
$zj.Range("D4").Value = 0.42

$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q4"
$zj.Range("C5").Value = 3
$zj.Range("D5").Value = 0.28

# ---------------------------------------------------------------------
# 2) Add the "2022-Q3" sheet right after "总计" by cloning the layout of
#    an existing quarter sheet (identical headers/styles), then
#    overwrite its data with the new quarter's numbers.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $zj)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("B2").Value = "008928"
$q3.Range("C2").Value = "泰达宏利中证主要消费红利指数A"
$q3.Range("D2").Value = "3.45"
$q3.Range("E2").Value = "93.07"
$q3.Range("F2").Value = "3.86"
$q3.Range("G2").Value = "0.1332"
$q3.Range("H2").Value = 9

$q3.Range("B3").Value = "008929"
$q3.Range("C3").Value = "泰达宏利中证主要消费红利指数C"
$q3.Range("D3").Value = "1.69"
$q3.Range("E3").Value = "93.07"
$q3.Range("F3").Value = "3.86"
$q3.Range("G3").Value = "0.0652"
$q3.Range("H3").Value = 9

$q3.Range("B4").Value = "501089"
$q3.Range("C4").Value = "方正富邦消费红利指数增强（LOF）"
$q3.Range("D4").Value = "0.31"
$q3.Range("E4").Value = "45.73"
$q3.Range("F4").Value = "1.87"
$q3.Range("G4").Value = "0.0058"
$q3.Range("H4").Value = 9
